$wb = $excel.ActiveWorkbook

# Sheet "展览": update "想去人数" (F column) counts for rows 2-5
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 256
$ws1.Range("F3").Value = 81
$ws1.Range("F4").Value = 865
$ws1.Range("F5").Value = 530

# Sheet "全部类型": same events appear here, but the music-concert row shifts
# row5 down, so the matching counts land on F2, F3, F4, F6
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 256
$ws4.Range("F3").Value = 81
$ws4.Range("F4").Value = 865
$ws4.Range("F6").Value = 530
